# "Generate Report for Handback" -- the localization-status report is
# regenerated after the de-de / zh-cn handback files were produced:
#   * Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   * The Latest Handback DateTime for each locale is refreshed
#   * The stale "handback file is not latest" Error Detail is cleared
#     now that the handback is in sync
#   * The Status / zh-cn / de-de columns widen to fit the new, longer
#     status text, and the (now empty) Error Detail column narrows

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet --------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

$ws1.Columns.Item(5).ColumnWidth = 29.15
$ws1.Columns.Item(6).ColumnWidth = 29.15

# ---- zh-cn sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = $newStatus
$ws2.Range("K2").Value = "2016-09-06 14:45:01"
$ws2.Range("P2").Value = ""

$ws2.Columns.Item(3).ColumnWidth = 29.15
$ws2.Columns.Item(16).ColumnWidth = 12.83

# ---- de-de sheet -------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = $newStatus
$ws3.Range("K2").Value = "2016-09-06 14:45:57"
$ws3.Range("P2").Value = ""

$ws3.Columns.Item(3).ColumnWidth = 29.15
$ws3.Columns.Item(16).ColumnWidth = 12.83
